$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting log entry for 3/1/2017 (row 8).
# Seed the row with the formatting (styles) of the prior entry (row 7),
# covering only columns A:E (column F / Additional Notes stays blank for
# this entry), then overwrite the values/content cell by cell.
$ws.Range("A7:E7").Copy($ws.Range("A8:E8"))

# Row is taller than prior rows because of the longer "Actions Required" text.
$ws.Rows("8:8").RowHeight = 105

# A8: Date
$ws.Range("A8").Value = 42795

# B8: Purpose of Meeting
$ws.Range("B8").Value = "Find out where we are on Research Discuss buying parts to begin prototyping Begin creating designs of circuits and software"

# C8: Attendees (same group as always)
$ws.Range("C8").Value = "Carlos, Courtnie, Lucas, Patrick"

# E8: Actions Required - rich text with bold names / bold+italic lead-in
$actionsText = "Please complete Research on at least one of the following by 3/1" + `
  ": Carlos" + `
  ": Smoke chamber, Direction/Location Algorithm Research & Docs, photoelectric sensor schematic; " + `
  "Courtnie" + `
  ": Battery/Power Monitoring, Alarm System Components Research, Docs and Schematics; " + `
  "Lucas: " + `
  "Processor Research & Docs; " + `
  "Patrick" + `
  ": Wireless Communication Research and Processor, Docs, and schematic;"

$ws.Range("E8").Value = $actionsText

$pos = 1

$len = ("Please complete Research on at least one of the following by 3/1").Length
$run = $ws.Range("E8").Characters($pos, $len)
$run.Font.Bold = $true
$run.Font.Italic = $true
$pos += $len

$len = (": Carlos").Length
$run = $ws.Range("E8").Characters($pos, $len)
$run.Font.Bold = $true
$pos += $len

$len = (": Smoke chamber, Direction/Location Algorithm Research & Docs, photoelectric sensor schematic; ").Length
$pos += $len

$len = ("Courtnie").Length
$run = $ws.Range("E8").Characters($pos, $len)
$run.Font.Bold = $true
$pos += $len

$len = (": Battery/Power Monitoring, Alarm System Components Research, Docs and Schematics; ").Length
$pos += $len

$len = ("Lucas: ").Length
$run = $ws.Range("E8").Characters($pos, $len)
$run.Font.Bold = $true
$pos += $len

$len = ("Processor Research & Docs; ").Length
$pos += $len

$len = ("Patrick").Length
$run = $ws.Range("E8").Characters($pos, $len)
$run.Font.Bold = $true
$pos += $len

$len = (": Wireless Communication Research and Processor, Docs, and schematic;").Length
$pos += $len

# D8: Resolution
$ws.Range("D8").Value = "Assigned Research and developent of schematics, will buy parts next week, requirements and specifications for software by next week"

# Selection / view state, as left after the meeting-log update.
$ws.Range("D8").Select()
$excel.ActiveWindow.ScrollRow = 6
